# "new migrations with tags"
# The authors/tags pivot tables no longer carry their own free-text
# "author" / "tag" columns -- they now reuse the common "name" column
# (matching how "collections" already stores its label in a "name"
# field). Update the two header-ish label cells accordingly; the
# now-unreferenced "author" / "tag" shared strings are dropped
# automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "name"
$ws.Range("F3").Value = "name"

# Move the active selection, matching where the author's cursor ended
# up after the edit.
[void]$ws.Range("G6").Select()
